$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 73. This shifts the existing
# rows 73-87 down to 74-88 (preserving their values/styles), matching the
# weekly update that prepends a new price report for this product.
$ws.Rows(73).Insert()

# Populate the newly inserted row 73 with this week's entry.
$ws.Cells.Item(73, 1).Value = 1
$ws.Cells.Item(73, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(73, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(73, 4).Value = 44722
$ws.Cells.Item(73, 5).Value = 15
$ws.Cells.Item(73, 6).Value = 100112038
$ws.Cells.Item(73, 7).Value = "Cebollín baby"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 300
$ws.Cells.Item(73, 11).Value = 2000
$ws.Cells.Item(73, 12).Value = 2500
$ws.Cells.Item(73, 13).Value = 2250
$ws.Cells.Item(73, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(73, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(73, 16).Value = 1125
$ws.Cells.Item(73, 17).Value = 2
$ws.Cells.Item(73, 18).Value = "Hortaliza"
